# Generate Report for Handback
# The 629d8d90-955e-4295-82e3-77be65549f95 file has been handed back and is
# in sync with en-US, so update its status + target/handback info on the
# per-language sheets (zh-cn, de-de) and reflect the same status on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$wsZh.Range("F2").Value = "629d8d90-955e-4295-82e3-77be65549f95.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7b019a96140972b8b7bff59e47f0d38f24242932/e2e/629d8d90-955e-4295-82e3-77be65549f95.md", "", "", "629d8d90-955e-4295-82e3-77be65549f95.md")

$wsZh.Range("G2").Value = "629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/46fe7b1a207394571f6566a107a0d60a762612b5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.zh-cn.xlf", "", "", "629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.zh-cn.xlf")

$wsZh.Range("H2").Value = "2016-03-19 14:31:53"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDe.Range("F2").Value = "629d8d90-955e-4295-82e3-77be65549f95.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7b019a96140972b8b7bff59e47f0d38f24242932/e2e/629d8d90-955e-4295-82e3-77be65549f95.md", "", "", "629d8d90-955e-4295-82e3-77be65549f95.md")

$wsDe.Range("G2").Value = "629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38de2dfea101b402eace4c64c4b0c99983d5e496/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.de-de.xlf", "", "", "629d8d90-955e-4295-82e3-77be65549f95.8060462a26539e58b8dcbc93cf2bd7e91913c09a.de-de.xlf")

$wsDe.Range("H2").Value = "2016-03-19 14:31:58"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
